# "Finished Spectra Calculation, Starting working on Summing all spectra"
#
# Adds an index column (E) running 1..100 alongside the per-frame centers
# (B) / sigmas (C) data in rows 2-101, then builds out a labeled summary
# block in rows 103-108: AVG / STD / Err AVG / (blank) / Disper / % on the
# left (column A, blue fill) mirrored by Width / STD / ARR Width / (blank)
# / Disper / % on the right (column D, green fill). The "Err AVG" formulas
# now divide by SQRT(E101) (the array width) instead of the literal 100,
# and a Dispersion-% formula is added for the sigmas column (C107).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: array index (1-based) for each data row 2..101 ---
for ($r = 2; $r -le 101; $r++) {
    $ws.Range("E$r").Value = $r - 1
}

# --- Left label column (A) ---
$ws.Range("A103").Value = "AVG"
$ws.Range("A104").Value = "STD"
$ws.Range("A105").Value = "Err AVG"
$ws.Range("A106").Value = ""
$ws.Range("A107").Value = "Disper"
$ws.Range("A108").Value = "%"

# --- Right label column (D) ---
$ws.Range("D103").Value = "Width"
$ws.Range("D104").Value = "STD"
$ws.Range("D105").Value = "ARR Width"
$ws.Range("D106").Value = ""
$ws.Range("D107").Value = "Disper"
$ws.Range("D108").Value = "%"

# --- Highlight the two label columns ---
$ws.Range("A103:A108").Interior.Color = 15773696   # FF00B0F0 (blue)
$ws.Range("D103:D108").Interior.Color = 5296274    # FF92D050 (green)

# --- Error-of-the-average formulas now reference the array width cell ---
$ws.Range("B105").Formula = "=B104/SQRT(E101)"
$ws.Range("C105").Formula = "=C104/SQRT(E101)"

# --- New dispersion-percent formula for the sigmas column ---
$ws.Range("C107").Formula = "=C104/C103*100"

# --- Leave the selection where the author left off ---
$ws.Range("C105").Select()
